# Updates cryptos list data (price + volume columns, and two row
# reorderings) per the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.930.12'
$ws.Range('E2').Value = '  -2.58%  '

# Row 3
$ws.Range('D3').Value = '2.918.15'
$ws.Range('E3').Value = '  -3.38%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.35'
$ws.Range('E5').Value = '  -1.70%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.04'
$ws.Range('E6').Value = '  -0.56%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.09%  '

# Row 8
$ws.Range('D8').Value = '2.918.64'
$ws.Range('E8').Value = '  -3.32%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.501'
$ws.Range('E9').Value = '  -3.30%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.72'
$ws.Range('E10').Value = '  +7.07%  '

# Row 11
$ws.Range('E11').Value = '  -3.39%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.447'
$ws.Range('E12').Value = '  -2.07%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000224'
$ws.Range('E13').Value = '  -3.78%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.47'
$ws.Range('E14').Value = '  -0.13%  '

# Row 15
$ws.Range('E15').Value = '  -0.37%  '

# Row 16
$ws.Range('D16').Value = '3.409.01'
$ws.Range('E16').Value = '  -3.16%  '

# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '61.074.59'
$ws.Range('E17').Value = '  -2.25%  '

# Row 18
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.84'
$ws.Range('E18').Value = '  -2.44%  '

# Row 19
$ws.Range('D19').Value = '2.925.67'
$ws.Range('E19').Value = '  -3.12%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '426.41'
$ws.Range('E20').Value = '  -5.32%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.69'
$ws.Range('E21').Value = '  -3.10%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.672'
$ws.Range('E22').Value = '  -2.23%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.18'
$ws.Range('E23').Value = '  -2.73%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.92'
$ws.Range('E24').Value = '  -1.32%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.99'
$ws.Range('E25').Value = '  -1.72%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.18'
$ws.Range('E26').Value = '  -4.43%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.88'

# Row 28
$ws.Range('E28').Value = '  -0.13%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.26'
$ws.Range('E29').Value = '  -0.94%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.25%  '

# Row 31
$ws.Range('E31').Value = '  -3.05%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.17'
$ws.Range('E32').Value = '  +2.96%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '26.77'
$ws.Range('E33').Value = '  -2.66%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.106'
$ws.Range('E34').Value = '  -3.94%  '

# Row 35
$ws.Range('D35').Value = '0.0₃0840'
$ws.Range('E35').Value = '  -1.54%  '

# Row 36
$ws.Range('E36').Value = '  -1.48%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.68'
$ws.Range('E37').Value = '  -2.72%  '

# Row 38
$ws.Range('E38').Value = '  +1.46%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '49.85'
$ws.Range('E39').Value = '  -1.12%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.03'
$ws.Range('E40').Value = '  -1.97%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.124'
$ws.Range('E41').Value = '  -0.22%  '

# Row 42
$ws.Range('E42').Value = '  -2.36%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.09'
$ws.Range('E43').Value = '  +1.96%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.288'
$ws.Range('E44').Value = '  +2.55%  '

# Row 45
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '374.28'
$ws.Range('E45').Value = '  -8.04%  '

# Row 46
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0347'
$ws.Range('E46').Value = '  -1.82%  '

# Row 47
$ws.Range('D47').Value = '2.657.42'
$ws.Range('E47').Value = '  -2.24%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.99'
$ws.Range('E48').Value = '  +0.04%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.45'
$ws.Range('E49').Value = '  +6.72%  '

# Row 50
$ws.Range('E50').Value = '  +0.01%  '

# Row 51
$ws.Range('E51').Value = '  -1.07%  '
